$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 108, shifting existing rows 108:208 down to 109:209
$ws.Rows.Item(108).Insert()

# Populate the newly inserted row 108 with the new data record
$ws.Cells.Item(108, 1).Value = 8
$ws.Cells.Item(108, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(108, 3).Value = "Coquimbo"
$ws.Cells.Item(108, 4).Value = 44484
$ws.Cells.Item(108, 5).Value = 4
$ws.Cells.Item(108, 6).Value = 100114013
$ws.Cells.Item(108, 7).Value = "Zanahoria"
$ws.Cells.Item(108, 8).Value = "Sin especificar"
$ws.Cells.Item(108, 9).Value = "Primera"
$ws.Cells.Item(108, 10).Value = 760
$ws.Cells.Item(108, 11).Value = 6500
$ws.Cells.Item(108, 12).Value = 7000
$ws.Cells.Item(108, 13).Value = 6750
$ws.Cells.Item(108, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(108, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(108, 16).Value = 338
$ws.Cells.Item(108, 17).Value = 20
$ws.Cells.Item(108, 18).Value = "Hortaliza"
